# Tau Dissertation Data - Participant 3 baseline trend correction
# Renames Sheet2, adds a new "3 corrected baseline" analysis sheet mirroring the
# layout of the "1 and 2 corrected baseline" sheet, and makes it the active tab.

$wb = $excel.ActiveWorkbook

# --- Rename Sheet2 ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "1 and 2 corrected baseline"

# --- Add the new sheet after it --------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "3 corrected baseline"

# --- Header row --------------------------------------------------------
$ws3.Range("A1").Value = "id"
$ws3.Range("B1").Value = "Label"
$ws3.Range("C1").Value = "S"
$ws3.Range("D1").Value = "PAIRS"
$ws3.Range("E1").Value = "TAU"
$ws3.Range("F1").Value = "TAUb"
$ws3.Range("G1").Value = "VARs"
$ws3.Range("H1").Value = "SD"
$ws3.Range("I1").Value = "SDtau"
$ws3.Range("J1").Value = "Z"
$ws3.Range("K1").Value = "P Value"
$ws3.Range("L1").Value = "CI 85%"
$ws3.Range("M1").Value = "CI 90%"

$ws3.Range("A2").Value = "trend:"

$ws3.Range("A3").Value = 0
$ws3.Range("C3").Value = 1
$ws3.Range("D3").Value = 3
$ws3.Range("E3").Value = 0.3333
$ws3.Range("F3").Value = 0.3333
$ws3.Range("G3").Value = 3.6667
$ws3.Range("H3").Value = 1.9149
$ws3.Range("I3").Value = 0.6383
$ws3.Range("J3").Value = 0.5222
$ws3.Range("K3").Value = 0.6015
$ws3.Range("L3").Value = "-0.586<>1"
$ws3.Range("M3").Value = "-0.717<>1"

$ws3.Range("A4").Value = 1
$ws3.Range("C4").Value = -1
$ws3.Range("D4").Value = 15
$ws3.Range("E4").Value = -0.0667
$ws3.Range("F4").Value = -0.1
$ws3.Range("G4").Value = 28.3333
$ws3.Range("H4").Value = 5.3229
$ws3.Range("I4").Value = 0.3549
$ws3.Range("J4").Value = -0.1879
$ws3.Range("K4").Value = 0.851
$ws3.Range("L4").Value = "-0.578<>0.444"
$ws3.Range("M4").Value = "-0.650<>0.517"

$ws3.Range("A5").Value = 2
$ws3.Range("C5").Value = 23
$ws3.Range("D5").Value = 28
$ws3.Range("E5").Value = 0.8214
$ws3.Range("F5").Value = 0.8364
$ws3.Range("G5").Value = 65.3333
$ws3.Range("H5").Value = 8.0829
$ws3.Range("I5").Value = 0.2887
$ws3.Range("J5").Value = 2.8455
$ws3.Range("K5").Value = 0.0044
$ws3.Range("L5").Value = "0.406<>1"
$ws3.Range("M5").Value = "0.347<>1"

$ws3.Range("A6").Value = 5
$ws3.Range("C6").Value = 23
$ws3.Range("D6").Value = 28
$ws3.Range("E6").Value = 0.8214
$ws3.Range("F6").Value = 0.8364
$ws3.Range("G6").Value = 65.3333
$ws3.Range("H6").Value = 8.0829
$ws3.Range("I6").Value = 0.2887
$ws3.Range("J6").Value = 2.8455
$ws3.Range("K6").Value = 0.0044
$ws3.Range("L6").Value = "0.406<>1"
$ws3.Range("M6").Value = "0.347<>1"

$ws3.Range("A7").Value = "phase:"

$ws3.Range("A8").Value = 3
$ws3.Range("C8").Value = 14
$ws3.Range("D8").Value = 30
$ws3.Range("E8").Value = 0.4667
$ws3.Range("F8").Value = 0.5
$ws3.Range("G8").Value = 140
$ws3.Range("H8").Value = 11.8322
$ws3.Range("I8").Value = 0.3944
$ws3.Range("J8").Value = 1.1832
$ws3.Range("K8").Value = 0.2367

$ws3.Range("A9").Value = 4
$ws3.Range("C9").Value = 42
$ws3.Range("D9").Value = 42
$ws3.Range("E9").Value = 1
$ws3.Range("F9").Value = 1
$ws3.Range("G9").Value = 196
$ws3.Range("H9").Value = 14
$ws3.Range("I9").Value = 0.3333
$ws3.Range("J9").Value = 3
$ws3.Range("K9").Value = 0.0027

$ws3.Range("A10").Value = "corrected baseline:"

$ws3.Range("A11").Value = 6
$ws3.Range("C11").Value = -10
$ws3.Range("D11").Value = 40
$ws3.Range("E11").Value = -0.25
$ws3.Range("F11").Value = -0.2597
$ws3.Range("G11").Value = 186.6667
$ws3.Range("H11").Value = 13.6626
$ws3.Range("I11").Value = 0.3416
$ws3.Range("J11").Value = -0.7319
$ws3.Range("K11").Value = 0.4642

$ws3.Range("A12").Value = "combined:"

# --- Shared strings must be created in the exact order the author typed
#     them so the workbook's shared-string table lines up: bottom summary
#     row first, then the corrected-baseline / phase rows, then the trend
#     rows - each row right-to-left. ------------------------------------

# Row 16 (right to left)
$ws3.Range("H16").Value = "0.0033<>0.8123"
$ws3.Range("G16").Value = "0.0683<>0.7473"
$ws3.Range("F16").Value = "0.1106<>0.7050"
$ws3.Range("A16").Value = "#3+#4+#6"

# Row 11 (right to left)
$ws3.Range("M11").Value = "-0.812<>0.312"
$ws3.Range("L11").Value = "-0.742<>0.242"
$ws3.Range("B11").Value = "P3 BL vs P3 I"

# Row 9 (right to left)
$ws3.Range("M9").Value = "0.452<>1"
$ws3.Range("L9").Value = "0.520<>1"
$ws3.Range("B9").Value = "P2 BL vs P2 I"

# Row 8 (right to left)
$ws3.Range("M8").Value = "-0.182<>1"
$ws3.Range("L8").Value = "-0.101<>1"
$ws3.Range("B8").Value = "P1 BL vs P1 I"

# Row 5/6 then 4 then 3
$ws3.Range("B5").Value = "P3 BL vs P3 BL"
$ws3.Range("B6").Value = "P3 BL vs P3 BL"
$ws3.Range("B4").Value = "P2 BL vs P2 BL"
$ws3.Range("B3").Value = "P1 BL vs P1 BL"

# --- Row 13: divider row of dashes -----------------------------------------
$ws3.Range("A13").Value = "-"
$ws3.Range("B13").Value = "-"
$ws3.Range("C13").Value = "-"
$ws3.Range("D13").Value = "-"
$ws3.Range("E13").Value = "-"
$ws3.Range("F13").Value = "-"
$ws3.Range("G13").Value = "-"
$ws3.Range("H13").Value = "-"
$ws3.Range("I13").Value = "-"
$ws3.Range("J13").Value = "-"
$ws3.Range("K13").Value = "-"
$ws3.Range("L13").Value = "-"
$ws3.Range("M13").Value = "-"

# --- Row 15: combined-summary header ----------------------------------------
$ws3.Range("A15").Value = "Label"
$ws3.Range("B15").Value = "Tau"
$ws3.Range("C15").Value = "Var-Tau"
$ws3.Range("D15").Value = "Z"
$ws3.Range("E15").Value = "P-Value"
$ws3.Range("F15").Value = "CI 85%"
$ws3.Range("G15").Value = "CI 90%"
$ws3.Range("H15").Value = "CI 95%"

# --- Row 16 numeric values ---------------------------------------------
$ws3.Range("B16").Value = 0.4078
$ws3.Range("C16").Value = 0.2064
$ws3.Range("D16").Value = 1.976
$ws3.Range("E16").Value = 0.0482

# --- Make the new sheet the active tab --------------------------------------
$ws3.Activate()
